$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-06-05 Monday", $false, $false, $false, $false, $false, $true, 1, $false, "2023-06-06 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("81-12=", $false, $false, $false, $false, $false, $true, 1, $false, "76-31=", 2) | Out-Null
$d.Content.Find.Execute("23+49=", $false, $false, $false, $false, $false, $true, 1, $false, "95-90=", 2) | Out-Null
$d.Content.Find.Execute("19+74=", $false, $false, $false, $false, $false, $true, 1, $false, "2+21=", 2) | Out-Null
$d.Content.Find.Execute("6+31=", $false, $false, $false, $false, $false, $true, 1, $false, "51+11=", 2) | Out-Null
$d.Content.Find.Execute("32+4=", $false, $false, $false, $false, $false, $true, 1, $false, "8+13=", 2) | Out-Null
$d.Content.Find.Execute("99-49=", $false, $false, $false, $false, $false, $true, 1, $false, "95-45=", 2) | Out-Null
$d.Content.Find.Execute("2+51=", $false, $false, $false, $false, $false, $true, 1, $false, "25-25=", 2) | Out-Null
$d.Content.Find.Execute("76+15=", $false, $false, $false, $false, $false, $true, 1, $false, "19+28=", 2) | Out-Null
$d.Content.Find.Execute("48+49=", $false, $false, $false, $false, $false, $true, 1, $false, "44-9=", 2) | Out-Null
$d.Content.Find.Execute("97-45=", $false, $false, $false, $false, $false, $true, 1, $false, "1+63=", 2) | Out-Null
$d.Content.Find.Execute("86+6=", $false, $false, $false, $false, $false, $true, 1, $false, "69-24=", 2) | Out-Null
$d.Content.Find.Execute("59-18=", $false, $false, $false, $false, $false, $true, 1, $false, "3+47=", 2) | Out-Null
$d.Content.Find.Execute("7+85=", $false, $false, $false, $false, $false, $true, 1, $false, "73-29=", 2) | Out-Null
$d.Content.Find.Execute("86-76=", $false, $false, $false, $false, $false, $true, 1, $false, "19+30=", 2) | Out-Null
$d.Content.Find.Execute("20+56=", $false, $false, $false, $false, $false, $true, 1, $false, "89-47=", 2) | Out-Null
$d.Content.Find.Execute("34+54=", $false, $false, $false, $false, $false, $true, 1, $false, "64-43=", 2) | Out-Null
$d.Content.Find.Execute("93-46=", $false, $false, $false, $false, $false, $true, 1, $false, "54-51=", 2) | Out-Null
$d.Content.Find.Execute("63+13=", $false, $false, $false, $false, $false, $true, 1, $false, "80-36=", 2) | Out-Null
$d.Content.Find.Execute("4+65=", $false, $false, $false, $false, $false, $true, 1, $false, "3+57=", 2) | Out-Null
$d.Content.Find.Execute("4+3=", $false, $false, $false, $false, $false, $true, 1, $false, "89-64=", 2) | Out-Null
$d.Content.Find.Execute("0+76=", $false, $false, $false, $false, $false, $true, 1, $false, "29+12=", 2) | Out-Null
$d.Content.Find.Execute("89-1=", $false, $false, $false, $false, $false, $true, 1, $false, "53-33=", 2) | Out-Null
$d.Content.Find.Execute("44+13=", $false, $false, $false, $false, $false, $true, 1, $false, "51+37=", 2) | Out-Null
$d.Content.Find.Execute("17+79=", $false, $false, $false, $false, $false, $true, 1, $false, "51-4=", 2) | Out-Null
$d.Content.Find.Execute("75+9=", $false, $false, $false, $false, $false, $true, 1, $false, "71-54=", 2) | Out-Null
$d.Content.Find.Execute("25+72=", $false, $false, $false, $false, $false, $true, 1, $false, "95-1=", 2) | Out-Null
$d.Content.Find.Execute("30+39=", $false, $false, $false, $false, $false, $true, 1, $false, "59-3=", 2) | Out-Null
$d.Content.Find.Execute("64+27=", $false, $false, $false, $false, $false, $true, 1, $false, "22+74=", 2) | Out-Null
$d.Content.Find.Execute("10+27=", $false, $false, $false, $false, $false, $true, 1, $false, "5+7=", 2) | Out-Null
$d.Content.Find.Execute("36+33=", $false, $false, $false, $false, $false, $true, 1, $false, "45-5=", 2) | Out-Null
$d.Content.Find.Execute("14+38=", $false, $false, $false, $false, $false, $true, 1, $false, "74+7=", 2) | Out-Null
$d.Content.Find.Execute("38+61=", $false, $false, $false, $false, $false, $true, 1, $false, "11+24=", 2) | Out-Null
$d.Content.Find.Execute("42+41=", $false, $false, $false, $false, $false, $true, 1, $false, "80-31=", 2) | Out-Null
$d.Content.Find.Execute("63-26=", $false, $false, $false, $false, $false, $true, 1, $false, "98-58=", 2) | Out-Null
$d.Content.Find.Execute("54+31=", $false, $false, $false, $false, $false, $true, 1, $false, "94-65=", 2) | Out-Null
$d.Content.Find.Execute("66-46=", $false, $false, $false, $false, $false, $true, 1, $false, "90-29=", 2) | Out-Null
$d.Content.Find.Execute("61-51=", $false, $false, $false, $false, $false, $true, 1, $false, "82-42=", 2) | Out-Null
$d.Content.Find.Execute("16-2=", $false, $false, $false, $false, $false, $true, 1, $false, "40-12=", 2) | Out-Null
$d.Content.Find.Execute("46-24=", $false, $false, $false, $false, $false, $true, 1, $false, "4+72=", 2) | Out-Null
$d.Content.Find.Execute("25+56=", $false, $false, $false, $false, $false, $true, 1, $false, "27+19=", 2) | Out-Null
$d.Content.Find.Execute("84-82=", $false, $false, $false, $false, $false, $true, 1, $false, "4+25=", 2) | Out-Null
$d.Content.Find.Execute("80+6=", $false, $false, $false, $false, $false, $true, 1, $false, "51+27=", 2) | Out-Null
$d.Content.Find.Execute("99-18=", $false, $false, $false, $false, $false, $true, 1, $false, "14+54=", 2) | Out-Null
$d.Content.Find.Execute("42+5=", $false, $false, $false, $false, $false, $true, 1, $false, "32+46=", 2) | Out-Null
$d.Content.Find.Execute("70+8=", $false, $false, $false, $false, $false, $true, 1, $false, "10+85=", 2) | Out-Null
$d.Content.Find.Execute("5+49=", $false, $false, $false, $false, $false, $true, 1, $false, "90-68=", 2) | Out-Null
$d.Content.Find.Execute("9+25=", $false, $false, $false, $false, $false, $true, 1, $false, "98-60=", 2) | Out-Null
$d.Content.Find.Execute("37+15=", $false, $false, $false, $false, $false, $true, 1, $false, "93-17=", 2) | Out-Null
$d.Content.Find.Execute("93-35=", $false, $false, $false, $false, $false, $true, 1, $false, "25+54=", 2) | Out-Null
$d.Content.Find.Execute("13+64=", $false, $false, $false, $false, $false, $true, 1, $false, "58-43=", 2) | Out-Null
$d.Content.Find.Execute("82-2=", $false, $false, $false, $false, $false, $true, 1, $false, "22+43=", 2) | Out-Null
$d.Content.Find.Execute("1+32=", $false, $false, $false, $false, $false, $true, 1, $false, "79-61=", 2) | Out-Null
$d.Content.Find.Execute("58-17=", $false, $false, $false, $false, $false, $true, 1, $false, "60+16=", 2) | Out-Null
$d.Content.Find.Execute("94-11=", $false, $false, $false, $false, $false, $true, 1, $false, "74-38=", 2) | Out-Null
$d.Content.Find.Execute("72+14=", $false, $false, $false, $false, $false, $true, 1, $false, "50+24=", 2) | Out-Null
$d.Content.Find.Execute("73-7=", $false, $false, $false, $false, $false, $true, 1, $false, "65-22=", 2) | Out-Null
$d.Content.Find.Execute("32+34=", $false, $false, $false, $false, $false, $true, 1, $false, "75-70=", 2) | Out-Null
$d.Content.Find.Execute("89-30=", $false, $false, $false, $false, $false, $true, 1, $false, "16+36=", 2) | Out-Null
$d.Content.Find.Execute("72+4=", $false, $false, $false, $false, $false, $true, 1, $false, "57+9=", 2) | Out-Null
$d.Content.Find.Execute("62+1=", $false, $false, $false, $false, $false, $true, 1, $false, "99-59=", 2) | Out-Null
$d.Content.Find.Execute("47+4=", $false, $false, $false, $false, $false, $true, 1, $false, "4-3=", 2) | Out-Null
$d.Content.Find.Execute("80-59=", $false, $false, $false, $false, $false, $true, 1, $false, "60-46=", 2) | Out-Null
$d.Content.Find.Execute("24+64=", $false, $false, $false, $false, $false, $true, 1, $false, "74+18=", 2) | Out-Null
$d.Content.Find.Execute("69-22=", $false, $false, $false, $false, $false, $true, 1, $false, "69+5=", 2) | Out-Null
$d.Content.Find.Execute("93-84=", $false, $false, $false, $false, $false, $true, 1, $false, "84-73=", 2) | Out-Null
$d.Content.Find.Execute("28+29=", $false, $false, $false, $false, $false, $true, 1, $false, "56-12=", 2) | Out-Null
$d.Content.Find.Execute("4+73=", $false, $false, $false, $false, $false, $true, 1, $false, "19+75=", 2) | Out-Null
$d.Content.Find.Execute("91-30=", $false, $false, $false, $false, $false, $true, 1, $false, "17-3=", 2) | Out-Null
$d.Content.Find.Execute("94-13=", $false, $false, $false, $false, $false, $true, 1, $false, "36+58=", 2) | Out-Null
$d.Content.Find.Execute("42+26=", $false, $false, $false, $false, $false, $true, 1, $false, "40+43=", 2) | Out-Null
$d.Content.Find.Execute("29+19=", $false, $false, $false, $false, $false, $true, 1, $false, "28+52=", 2) | Out-Null
$d.Content.Find.Execute("30+54=", $false, $false, $false, $false, $false, $true, 1, $false, "93-55=", 2) | Out-Null
$d.Content.Find.Execute("97-86=", $false, $false, $false, $false, $false, $true, 1, $false, "52+26=", 2) | Out-Null
$d.Content.Find.Execute("1+74=", $false, $false, $false, $false, $false, $true, 1, $false, "61-59=", 2) | Out-Null
$d.Content.Find.Execute("82+3=", $false, $false, $false, $false, $false, $true, 1, $false, "74-20=", 2) | Out-Null
$d.Content.Find.Execute("39-6=", $false, $false, $false, $false, $false, $true, 1, $false, "98-27=", 2) | Out-Null
$d.Content.Find.Execute("39+10=", $false, $false, $false, $false, $false, $true, 1, $false, "27+41=", 2) | Out-Null
$d.Content.Find.Execute("95-40=", $false, $false, $false, $false, $false, $true, 1, $false, "0+61=", 2) | Out-Null
$d.Content.Find.Execute("37+2=", $false, $false, $false, $false, $false, $true, 1, $false, "11+79=", 2) | Out-Null
$d.Content.Find.Execute("43-34=", $false, $false, $false, $false, $false, $true, 1, $false, "70-13=", 2) | Out-Null
$d.Content.Find.Execute("9+24=", $false, $false, $false, $false, $false, $true, 1, $false, "11-10=", 2) | Out-Null
$d.Content.Find.Execute("16+17=", $false, $false, $false, $false, $false, $true, 1, $false, "85-20=", 2) | Out-Null
$d.Content.Find.Execute("0+71=", $false, $false, $false, $false, $false, $true, 1, $false, "40+46=", 2) | Out-Null
$d.Content.Find.Execute("21-12=", $false, $false, $false, $false, $false, $true, 1, $false, "51-34=", 2) | Out-Null
$d.Content.Find.Execute("33-31=", $false, $false, $false, $false, $false, $true, 1, $false, "48-18=", 2) | Out-Null
$d.Content.Find.Execute("10-9=", $false, $false, $false, $false, $false, $true, 1, $false, "39+11=", 2) | Out-Null
$d.Content.Find.Execute("84-11=", $false, $false, $false, $false, $false, $true, 1, $false, "25+14=", 2) | Out-Null
$d.Content.Find.Execute("60+12=", $false, $false, $false, $false, $false, $true, 1, $false, "32+5=", 2) | Out-Null
$d.Content.Find.Execute("34+5=", $false, $false, $false, $false, $false, $true, 1, $false, "76-22=", 2) | Out-Null
$d.Content.Find.Execute("99-16=", $false, $false, $false, $false, $false, $true, 1, $false, "95+3=", 2) | Out-Null
$d.Content.Find.Execute("92-74=", $false, $false, $false, $false, $false, $true, 1, $false, "57+12=", 2) | Out-Null
$d.Content.Find.Execute("95-21=", $false, $false, $false, $false, $false, $true, 1, $false, "89-33=", 2) | Out-Null
$d.Content.Find.Execute("67-5=", $false, $false, $false, $false, $false, $true, 1, $false, "83-60=", 2) | Out-Null
$d.Content.Find.Execute("25+41=", $false, $false, $false, $false, $false, $true, 1, $false, "13-5=", 2) | Out-Null
$d.Content.Find.Execute("23+27=", $false, $false, $false, $false, $false, $true, 1, $false, "50-20=", 2) | Out-Null
$d.Content.Find.Execute("19+1=", $false, $false, $false, $false, $false, $true, 1, $false, "40+13=", 2) | Out-Null
$d.Content.Find.Execute("47+19=", $false, $false, $false, $false, $false, $true, 1, $false, "28+67=", 2) | Out-Null
$d.Content.Find.Execute("50-48=", $false, $false, $false, $false, $false, $true, 1, $false, "91-16=", 2) | Out-Null
$d.Content.Find.Execute("50-33=", $false, $false, $false, $false, $false, $true, 1, $false, "72-61=", 2) | Out-Null
$d.Content.Find.Execute("54-49=", $false, $false, $false, $false, $false, $true, 1, $false, "81-39=", 2) | Out-Null
